$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 301.29092
$ws.Range("I8").Value = 341.33334
$ws.Range("J8").Value = 298.98077
$ws.Range("K8").Value = 1024.00002
$ws.Range("L8").Value = 896.94231
$ws.Range("M8").Value = -885.0000199999999
$ws.Range("N8").Value = -1174.94231
$ws.Range("H94").Value = 3599.6155
$ws.Range("I94").Value = 1399.5834
$ws.Range("K94").Value = 1399.5834
$ws.Range("M94").Value = -948.5834
$ws.Range("H103").Value = 797.7143
$ws.Range("I103").Value = 570.8333
$ws.Range("K103").Value = 1712.4999
$ws.Range("M103").Value = -1126.4999
$ws.Range("H106").Value = 1636.8
$ws.Range("I106").Value = 1636.8
$ws.Range("K106").Value = 1636.8
$ws.Range("M106").Value = -1005.8
$ws.Range("H116").Value = 23752.857
$ws.Range("I116").Value = 25717.785
$ws.Range("J116").Value = 19823
$ws.Range("K116").Value = 25717.785
$ws.Range("L116").Value = 19823
$ws.Range("M116").Value = -22275.785
$ws.Range("N116").Value = -26707
$ws.Range("H132").Value = 25389.8
$ws.Range("I132").Value = 28027.592
$ws.Range("J132").Value = 1649.6666
$ws.Range("K132").Value = 84082.776
$ws.Range("L132").Value = 4948.9998
$ws.Range("M132").Value = -81552.776
$ws.Range("N132").Value = -10008.9998
$ws.Range("H137").Value = 23250.5
$ws.Range("I137").Value = 38201.625
$ws.Range("K137").Value = 114604.875
$ws.Range("M137").Value = -112054.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38587.395
$ws.Range("I32").Value = 42978
$ws.Range("J32").Value = 1999
$ws.Range("K32").Value = 42978
$ws.Range("L32").Value = 1999
$ws.Range("M32").Value = -42691
$ws.Range("N32").Value = -2573
$ws.Range("H74").Value = 218797.83
$ws.Range("I74").Value = 286249.84
$ws.Range("K74").Value = 286249.84
$ws.Range("M74").Value = -285375.84
$ws.Range("H77").Value = 218797.83
$ws.Range("I77").Value = 286249.84
$ws.Range("K77").Value = 1431249.2
$ws.Range("M77").Value = -1426881.2
$ws.Range("H102").Value = 2332.16
$ws.Range("I102").Value = 2318.5908
$ws.Range("K102").Value = 2318.5908
$ws.Range("M102").Value = -696.5907999999999
$ws.Range("H132").Value = 1296.9474
$ws.Range("I132").Value = 1036.6857
$ws.Range("J132").Value = 4333.3335
$ws.Range("K132").Value = 3110.0571
$ws.Range("L132").Value = 13000.0005
$ws.Range("M132").Value = -580.0571
$ws.Range("N132").Value = -18060.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 19369.055
$ws.Range("I20").Value = 33056.2
$ws.Range("K20").Value = 33056.2
$ws.Range("M20").Value = -32809.2
$ws.Range("H99").Value = 2799.0588
$ws.Range("I99").Value = 1524.8182
$ws.Range("J99").Value = 5135.1665
$ws.Range("K99").Value = 1524.8182
$ws.Range("L99").Value = 5135.1665
$ws.Range("M99").Value = -26.81819999999993
$ws.Range("N99").Value = -8131.1665
$ws.Range("H105").Value = 2225.926
$ws.Range("I105").Value = 1424.9445
$ws.Range("K105").Value = 1424.9445
$ws.Range("M105").Value = 322.0554999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4001202.5
$ws.Range("I31").Value = 4546400.5
$ws.Range("K31").Value = 4546400.5
$ws.Range("M31").Value = -4546105.5
$ws.Range("H34").Value = 4001202.5
$ws.Range("I34").Value = 4546400.5
$ws.Range("K34").Value = 4546400.5
$ws.Range("M34").Value = -4546198.5
$ws.Range("H58").Value = 15861.131
$ws.Range("I58").Value = 1490.4736
$ws.Range("K58").Value = 1490.4736
$ws.Range("M58").Value = -1287.4736
$ws.Range("H86").Value = 34852.37
$ws.Range("I86").Value = 51104.465
$ws.Range("J86").Value = 14537.25
$ws.Range("K86").Value = 51104.465
$ws.Range("L86").Value = 14537.25
$ws.Range("M86").Value = -49981.465
$ws.Range("N86").Value = -16783.25
$ws.Range("H89").Value = 34852.37
$ws.Range("I89").Value = 51104.465
$ws.Range("J89").Value = 14537.25
$ws.Range("K89").Value = 255522.325
$ws.Range("L89").Value = 72686.25
$ws.Range("M89").Value = -249906.325
$ws.Range("N89").Value = -83918.25
$ws.Range("H107").Value = 838.2353000000001
$ws.Range("I107").Value = 723.9091
$ws.Range("J107").Value = 1047.8334
$ws.Range("K107").Value = 723.9091
$ws.Range("L107").Value = 1047.8334
$ws.Range("M107").Value = 1196.0909
$ws.Range("N107").Value = -4887.8334
$ws.Range("H132").Value = 51769.55
$ws.Range("I132").Value = 68026.13
$ws.Range("K132").Value = 204078.39
$ws.Range("M132").Value = -201548.39
$ws.Range("H136").Value = 15861.131
$ws.Range("I136").Value = 1490.4736
$ws.Range("K136").Value = 4471.4208
$ws.Range("M136").Value = -1921.4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1005.5
$ws.Range("I50").Value = 1005
$ws.Range("K50").Value = 3015
$ws.Range("M50").Value = -2534
$ws.Range("H53").Value = 1005.5
$ws.Range("I53").Value = 1005
$ws.Range("K53").Value = 3015
$ws.Range("M53").Value = -2534
$ws.Range("H62").Value = 2975
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2975
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 750
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H114").Value = 756.1111
$ws.Range("I114").Value = 200
$ws.Range("J114").Value = 1201
$ws.Range("K114").Value = 600
$ws.Range("L114").Value = 3603
$ws.Range("M114").Value = 2654
$ws.Range("N114").Value = -10111
$ws.Range("H131").Value = 2765.524
$ws.Range("I131").Value = 2845.875
$ws.Range("J131").Value = 2508.4
$ws.Range("K131").Value = 8537.625
$ws.Range("L131").Value = 7525.200000000001
$ws.Range("M131").Value = -3497.625
$ws.Range("N131").Value = -17605.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 14500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1537.3334
$ws.Range("I55").Value = 298.6
$ws.Range("J55").Value = 3085.75
$ws.Range("K55").Value = 298.6
$ws.Range("L55").Value = 3085.75
$ws.Range("M55").Value = -125.6
$ws.Range("N55").Value = -3431.75
$ws.Range("H136").Value = 4017.7856
$ws.Range("I136").Value = 3717.3635
$ws.Range("J136").Value = 5119.3335
$ws.Range("K136").Value = 11152.0905
$ws.Range("L136").Value = 15358.0005
$ws.Range("M136").Value = -8602.0905
$ws.Range("N136").Value = -20458.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 32076.25
$ws.Range("I136").Value = 43788.355
$ws.Range("K136").Value = 131365.065
$ws.Range("M136").Value = -128815.065
